$d = $word.ActiveDocument
$para = $d.Paragraphs(1).Range
$end = $para.End - 1
$r = $d.Range($end, $end)
$r.InsertAfter(" ")
